$d = $word.ActiveDocument
$start = 339
$r = $d.Range($start, $start+1)
$r.Text = "th"

$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute(" of May, 2015", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$spaceRng = $d.Range($rng.Start, $rng.Start + 1)
$restRng = $d.Range($rng.Start + 1, $rng.End)
$restRng.Bold = 1
$restRng.Bold = 0
